$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06729733333333333
$ws.Range("H2").Value = 0.201892
$ws.Range("I2").Value = 0.01373511018321553
$ws.Range("J2").Value = 0.01373511018321553
$ws.Range("M2").Value = 1.952294
$ws.Range("N2").Value = 5.856882
$ws.Range("O2").Value = 0.07575070565202183
$ws.Range("P2").Value = 0.07575070565202184
$ws.Range("Q2").Value = 0.1313841800826667
$ws.Range("R2").Value = 1.182457620744
$ws.Range("S2").Value = 0.001040444288586848
$ws.Range("T2").Value = 0.001040444288586848

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06729733333333333
$ws.Range("H3").Value = 0.201892
$ws.Range("I3").Value = 0.01373511018321553
$ws.Range("J3").Value = 0.01373511018321553
$ws.Range("O3").Value = 0.5679402069281436
$ws.Range("P3").Value = 0.5679402069281437
$ws.Range("Q3").Value = 0.9850516609840001
$ws.Range("R3").Value = 8.865464948856
$ws.Range("S3").Value = 0.007800721319636283
$ws.Range("T3").Value = 0.007800721319636283

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06729733333333333
$ws.Range("H4").Value = 0.201892
$ws.Range("I4").Value = 0.01373511018321553
$ws.Range("J4").Value = 0.01373511018321553
$ws.Range("M4").Value = 9.009963
$ws.Range("N4").Value = 27.029889
$ws.Range("O4").Value = 0.3495944028658632
$ws.Range("P4").Value = 0.3495944028658634
$ws.Range("Q4").Value = 0.606346483332
$ws.Range("R4").Value = 5.457118349988001
$ws.Range("S4").Value = 0.004801717642798071
$ws.Range("T4").Value = 0.004801717642798073

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("E5").Value = 1.0
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06729733333333333
$ws.Range("H5").Value = 0.201892
$ws.Range("I5").Value = 0.01373511018321553
$ws.Range("J5").Value = 0.01373511018321553
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.173055
$ws.Range("N5").Value = 0.519165
$ws.Range("O5").Value = 0.006714684553971194
$ws.Range("P5").Value = 0.006714684553971196
$ws.Range("Q5").Value = 0.01164614002
$ws.Range("R5").Value = 0.10481526018
$ws.Range("S5").Value = 0.00009222693219432979
$ws.Range("T5").Value = 0.00009222693219432981

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 4.83236
$ws.Range("H6").Value = 14.49708
$ws.Range("I6").Value = 0.9862648898167845
$ws.Range("J6").Value = 0.9862648898167844
$ws.Range("M6").Value = 1.952294
$ws.Range("N6").Value = 5.856882
$ws.Range("O6").Value = 0.07575070565202183
$ws.Range("P6").Value = 0.07575070565202184
$ws.Range("Q6").Value = 9.43418743384
$ws.Range("R6").Value = 84.90768690456
$ws.Range("S6").Value = 0.07471026136343498
$ws.Range("T6").Value = 0.07471026136343499

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 4.83236
$ws.Range("H7").Value = 14.49708
$ws.Range("I7").Value = 0.9862648898167845
$ws.Range("J7").Value = 0.9862648898167844
$ws.Range("O7").Value = 0.5679402069281436
$ws.Range("P7").Value = 0.5679402069281437
$ws.Range("Q7").Value = 70.73273202216001
$ws.Range("R7").Value = 636.5945881994401
$ws.Range("S7").Value = 0.5601394856085073
$ws.Range("T7").Value = 0.5601394856085075

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 4.83236
$ws.Range("H8").Value = 14.49708
$ws.Range("I8").Value = 0.9862648898167845
$ws.Range("J8").Value = 0.9862648898167844
$ws.Range("M8").Value = 9.009963
$ws.Range("N8").Value = 27.029889
$ws.Range("O8").Value = 0.3495944028658632
$ws.Range("P8").Value = 0.3495944028658634
$ws.Range("Q8").Value = 43.53938480268001
$ws.Range("R8").Value = 391.8544632241201
$ws.Range("S8").Value = 0.3447926852230652
$ws.Range("T8").Value = 0.3447926852230653

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 4.83236
$ws.Range("H9").Value = 14.49708
$ws.Range("I9").Value = 0.9862648898167845
$ws.Range("J9").Value = 0.9862648898167844
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.173055
$ws.Range("N9").Value = 0.519165
$ws.Range("O9").Value = 0.006714684553971194
$ws.Range("P9").Value = 0.006714684553971196
$ws.Range("Q9").Value = 0.8362640598
$ws.Range("R9").Value = 7.5263765382
$ws.Range("S9").Value = 0.006622457621776864
$ws.Range("T9").Value = 0.006622457621776865
